$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- Fill in the real-effort / estimate numbers and the "done" flag (col F) ---
$ws.Range("F16").Value = 1
$ws.Range("F18").Value = 1

$ws.Range("C19").Value = 2.5
$ws.Range("D19").Value = 2
$ws.Range("F19").Value = 1

$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 1.5
$ws.Range("F20").Value = 1

$ws.Range("C21").Value = 1.5
$ws.Range("F21").Value = 1

$ws.Range("C22").Value = 1.5
$ws.Range("F22").Value = 1

$ws.Range("C23").Value = 2.5
$ws.Range("D23").Value = 2
$ws.Range("F23").Value = 1

$ws.Range("C25").Value = 1
$ws.Range("D25").Value = 1.5
$ws.Range("F25").Value = 1

$ws.Range("C26").Value = 2.5
$ws.Range("D26").Value = 0.5
$ws.Range("F26").Value = 1

$ws.Range("C27").Value = 2.5
$ws.Range("D27").Value = 0.5
$ws.Range("F27").Value = 1

$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 0.25
$ws.Range("F28").Value = 1

# --- Update the saved view state (scroll position + selection) ---
$ws.Range("C6:G6").Select()
$ws.Application.ActiveWindow.ScrollRow = 10
